# Auto-generated edit script: update cryptocurrency price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.472.43'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '3.082.44'
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''394.16'
$ws.Range("E5").Value = '  +2.43%  '

$ws.Range("D6").Value = '''102.68'
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").Value = '''0.536'
$ws.Range("E7").Value = '  -1.69%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '''0.589'
$ws.Range("E9").Value = '  +0.46%  '

$ws.Range("D10").Value = '''37.59'
$ws.Range("E10").Value = '  +1.58%  '

$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").Value = '''0.0854'
$ws.Range("E12").Value = '  -1.49%  '

$ws.Range("D13").Value = '3.556.03'
$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("E14").Value = '  -0.83%  '

$ws.Range("D15").Value = '''7.71'
$ws.Range("E15").Value = '  -0.57%  '

$ws.Range("E16").Value = '  +4.94%  '

$ws.Range("D17").Value = '3.047.23'
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").Value = '''10.59'
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").Value = '51.470.22'
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").Value = '''3.18'
$ws.Range("E20").Value = '  +1.91%  '

$ws.Range("D21").Value = '''12.42'
$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("D22").Value = '0.0₃0961'
$ws.Range("E22").Value = '  -0.53%  '

$ws.Range("D23").Value = '''70.31'
$ws.Range("E23").Value = '  +0.29%  '

$ws.Range("D24").Value = '''265.06'
$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("D26").Value = '''7.90'
$ws.Range("E26").Value = '  -6.88%  '

$ws.Range("D27").Value = '''27.01'
$ws.Range("E27").Value = '  +1.80%  '

$ws.Range("D28").Value = '''7.19'
$ws.Range("E28").Value = '  -2.18%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("E30").Value = '  -4.38%  '

$ws.Range("D31").Value = '''10.73'
$ws.Range("E31").Value = '  +4.11%  '

$ws.Range("E32").Value = '  -2.87%  '

$ws.Range("D33").Value = '''0.0495'
$ws.Range("E33").Value = '  +11.48%  '

$ws.Range("D34").Value = '''36.48'
$ws.Range("E34").Value = '  +6.50%  '

$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("D36").Value = '''49.90'
$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").Value = '''3.34'
$ws.Range("E38").Value = '  -1.16%  '

$ws.Range("D39").Value = '''4.03'
$ws.Range("E39").Value = '  +9.24%  '

$ws.Range("D40").Value = '''0.289'
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("D41").Value = '''129.29'
$ws.Range("E41").Value = '  +0.97%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '''1.85'
$ws.Range("E42").Value = '  -1.32%  '

$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '''16.71'
$ws.Range("E43").Value = '  -2.56%  '

$ws.Range("E44").Value = '  -0.75%  '

$ws.Range("D45").Value = '''2.54'
$ws.Range("E45").Value = '  -0.16%  '

$ws.Range("D46").Value = '''21.80'
$ws.Range("E46").Value = '  -0.71%  '

$ws.Range("D47").Value = '''2.53'
$ws.Range("E47").Value = '  +0.36%  '

$ws.Range("D48").Value = '''2.05'
$ws.Range("E48").Value = '  -1.75%  '

$ws.Range("D49").Value = '2.072.70'
$ws.Range("E49").Value = '  +1.68%  '

$ws.Range("D50").Value = '''0.0539'
$ws.Range("E50").Value = '  +37.87%  '

$ws.Range("D51").Value = '''0.902'
$ws.Range("E51").Value = '  +9.69%  '
